$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "64.893.76"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.556.33"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.86%  "

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "600.77"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.46%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "135.02"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.33%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.561.03"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.13%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.494"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.55%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.123"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.01%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.387"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.26%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.157.10"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.76%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000182"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.51%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.557.69"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.82%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "27.02"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.10%  "

# Row 17
$ws.Range("E17").Value = "  +0.80%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "64.627.54"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.07"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.78%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.40"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +5.19%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.83"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "386.26"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.578"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.87%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.697.10"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "74.30"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.84%  "

# Row 26
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0000117"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +9.77%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +5.34%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.30"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +4.37%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.39"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.39%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.47"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +23.47%  "

# Row 33
$ws.Range("B33").Value = "RenzoRestakedETH"
$ws.Range("C33").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.561.43"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.49%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "24.00"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.74%  "

# Row 35
$ws.Range("E35").Value = "  -0.01%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.144"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.76%  "

# Row 37
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.94"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.94%  "

# Row 38
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "169.39"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.55%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.54"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.25%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.99"
$c.Style = "Normal"

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0805"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.95%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "27.24"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +17.12%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.827"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.55%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "42.66"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.51%  "

# Row 45
$ws.Range("E45").Value = "  -0.05%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.46"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.35%  "

# Row 47
$ws.Range("E47").Value = "  +7.34%  "

# Row 48
$ws.Range("E48").Value = "  +1.99%  "

# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.469.87"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +11.43%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "6.93"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +5.34%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.38"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +11.82%  "
